$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells stay formatted as text, matching the source data
# which stores prices as literal strings (e.g. "1.010", "20.611.58").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '20.636.04'
$ws.Range("E2").Value = '  +2.70%  '
$ws.Range("D3").Value = '1.471.64'
$ws.Range("E3").Value = '  +2.99%  '
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.79%  '
$ws.Range("D5").Value = '0.9539'
$ws.Range("E5").Value = '  -4.69%  '
$ws.Range("D6").Value = '282.13'
$ws.Range("E6").Value = '  +2.70%  '
$ws.Range("D7").Value = '0.3725'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("D8").Value = '0.3200'
$ws.Range("E8").Value = '  +3.53%  '
$ws.Range("D9").Value = '41.91'
$ws.Range("E9").Value = '  +4.19%  '
$ws.Range("E10").Value = '  +5.26%  '
$ws.Range("D11").Value = '0.06724'
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Value = '5.649'
$ws.Range("E13").Value = '  +4.72%  '
$ws.Range("E14").Value = '  +7.00%  '
$ws.Range("E15").Value = '  +1.89%  '
$ws.Range("D16").Value = '1.478.36'
$ws.Range("E16").Value = '  +3.38%  '
$ws.Range("E17").Value = '  +3.02%  '
$ws.Range("D18").Value = '0.05789'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").Value = '73.01'
$ws.Range("E19").Value = '  -3.41%  '
$ws.Range("D20").Value = '0.9561'
$ws.Range("E20").Value = '  -4.45%  '
$ws.Range("D21").Value = '5.729'
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("E22").Value = '  +2.68%  '
$ws.Range("D23").Value = '11.22'
$ws.Range("E23").Value = '  +0.89%  '
$ws.Range("D24").Value = '2.297'
$ws.Range("E24").Value = '  -1.70%  '
$ws.Range("D25").Value = '20.734.52'
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("D26").Value = '2.334'
$ws.Range("E26").Value = '  +2.18%  '
$ws.Range("D27").Value = '138.01'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("E28").Value = '  +4.44%  '
$ws.Range("D29").Value = '1.640.84'
$ws.Range("E29").Value = '  +3.09%  '
$ws.Range("D30").Value = '114.05'
$ws.Range("E30").Value = '  +4.29%  '
$ws.Range("D31").Value = '3.990'
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").Value = '5.385'
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("D33").Value = '0.8462'
$ws.Range("E33").Value = '  -6.94%  '
$ws.Range("D34").Value = '1.654'
$ws.Range("E34").Value = '  +27.33%  '
$ws.Range("D35").Value = '0.07879'
$ws.Range("E35").Value = '  +1.32%  '
$ws.Range("D36").Value = '0.06116'
$ws.Range("E36").Value = '  +7.39%  '
$ws.Range("D37").Value = '4.964'
$ws.Range("E37").Value = '  +4.16%  '
$ws.Range("D38").Value = '10.80'
$ws.Range("E38").Value = '  -5.18%  '
$ws.Range("D39").Value = '0.02084'
$ws.Range("E39").Value = '  +2.65%  '
$ws.Range("D40").Value = '1.131'
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("D41").Value = '0.9639'
$ws.Range("E41").Value = '  -3.63%  '
$ws.Range("D42").Value = '0.1914'
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("D43").Value = '7.418'
$ws.Range("E43").Value = '  -11.80%  '
$ws.Range("D44").Value = '0.5469'
$ws.Range("E44").Value = '  +2.46%  '
$ws.Range("D45").Value = '12.56'
$ws.Range("E45").Value = '  +2.85%  '
$ws.Range("D46").Value = '3.600'
$ws.Range("E46").Value = '  +1.58%  '
$ws.Range("D47").Value = '121.82'
$ws.Range("E47").Value = '  +11.12%  '
$ws.Range("D48").Value = '0.5398'
$ws.Range("E48").Value = '  +4.96%  '
$ws.Range("D49").Value = '1.843'
$ws.Range("E49").Value = '  +3.67%  '
$ws.Range("D50").Value = '0.06463'
$ws.Range("E50").Value = '  +4.15%  '
$ws.Range("E51").Value = '  +0.49%  '
